$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Merge the split runs in the two TODO bullet paragraphs so each
#    paragraph holds a single run of text (matches the authored diff).
# ------------------------------------------------------------------

$mergeTargets = @(
    "* root with conjuncts `u{2013} OK",
    "* relative clauses - working"
)

foreach ($text in $mergeTargets) {
    # Re-write the text (Find/Replace merges any runs it spans into one
    # run); the replacement collapses formatting so immediately after we
    # flip Bold on/off on the same text to coax the engine into emitting
    # an (empty) rPr element on the resulting run, matching the source
    # paragraphs which all carry an explicit (empty) <w:rPr/>.
    $d.Content.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, $text, 2) | Out-Null

    $rng = $d.Content
    $rng.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $rng.Bold = 1
    $rng.Bold = 0
}

# ------------------------------------------------------------------
# 2) Normal style: stop allowing punctuation to overflow the margin
#    (w:overflowPunct false).
# ------------------------------------------------------------------
$normal = $d.Styles("Normal")
$normal.ParagraphFormat.HangingPunctuation = 0

# ------------------------------------------------------------------
# 3) Mint 20 new (empty) character styles, ListLabel 161 .. ListLabel 180,
#    mirroring the ListLabel1..160 family already in the template.
# ------------------------------------------------------------------
for ($i = 161; $i -le 180; $i++) {
    $style = $d.Styles.Add("ListLabel $i", 2)
    $style.QuickStyle = $true
}
